$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30
$ws.Range("A30").Value = 10
$ws.Range("B30").Value = 'Vega Modelo de Temuco'
$ws.Range("C30").Value = 'La Araucanía'
$ws.Range("D30").Value = 44606
$ws.Range("E30").Value = 9
$ws.Range("F30").Value = 100114002
$ws.Range("G30").Value = 'Camote'
$ws.Range("H30").Value = 'Sin especificar'
$ws.Range("I30").Value = 'Primera'
$ws.Range("J30").Value = 50
$ws.Range("K30").Value = 18000
$ws.Range("L30").Value = 18000
$ws.Range("M30").Value = 18000
$ws.Range("N30").Value = '$/malla 20 kilos'
$ws.Range("O30").Value = 'Perú'
$ws.Range("P30").Value = 900
$ws.Range("Q30").Value = 20
$ws.Range("R30").Value = 'Hortaliza'

# Row 31
$ws.Range("A31").Value = 10
$ws.Range("B31").Value = 'Vega Modelo de Temuco'
$ws.Range("C31").Value = 'La Araucanía'
$ws.Range("D31").Value = 44188
$ws.Range("E31").Value = 9
$ws.Range("F31").Value = 100114002
$ws.Range("G31").Value = 'Camote'
$ws.Range("H31").Value = 'Sin especificar'
$ws.Range("I31").Value = 'Primera'
$ws.Range("J31").Value = 20
$ws.Range("K31").Value = 20000
$ws.Range("L31").Value = 20000
$ws.Range("M31").Value = 20000
$ws.Range("N31").Value = '$/caja 15 kilos granel'
$ws.Range("O31").Value = 'Región de Arica y Parinacota'
$ws.Range("P31").Value = 1333
$ws.Range("Q31").Value = 15
$ws.Range("R31").Value = 'Hortaliza'

# Row 32
$ws.Range("A32").Value = 10
$ws.Range("B32").Value = 'Vega Modelo de Temuco'
$ws.Range("C32").Value = 'La Araucanía'
$ws.Range("D32").Value = 44175
$ws.Range("E32").Value = 9
$ws.Range("F32").Value = 100114002
$ws.Range("G32").Value = 'Camote'
$ws.Range("H32").Value = 'Sin especificar'
$ws.Range("I32").Value = 'Primera'
$ws.Range("J32").Value = 20
$ws.Range("K32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = 20000
$ws.Range("N32").Value = '$/caja 15 kilos granel'
$ws.Range("O32").Value = 'Región de Arica y Parinacota'
$ws.Range("P32").Value = 1333
$ws.Range("Q32").Value = 15
$ws.Range("R32").Value = 'Hortaliza'

# Row 33
$ws.Range("A33").Value = 10
$ws.Range("B33").Value = 'Vega Modelo de Temuco'
$ws.Range("C33").Value = 'La Araucanía'
$ws.Range("D33").Value = 44579
$ws.Range("E33").Value = 9
$ws.Range("F33").Value = 100114002
$ws.Range("G33").Value = 'Camote'
$ws.Range("H33").Value = 'Sin especificar'
$ws.Range("I33").Value = 'Primera'
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 20000
$ws.Range("L33").Value = 20000
$ws.Range("M33").Value = 20000
$ws.Range("N33").Value = '$/malla 20 kilos'
$ws.Range("O33").Value = 'Perú'
$ws.Range("P33").Value = 1000
$ws.Range("Q33").Value = 20
$ws.Range("R33").Value = 'Hortaliza'

# Row 34
$ws.Range("A34").Value = 10
$ws.Range("B34").Value = 'Vega Modelo de Temuco'
$ws.Range("C34").Value = 'La Araucanía'
$ws.Range("D34").Value = 44438
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = 100114002
$ws.Range("G34").Value = 'Camote'
$ws.Range("H34").Value = 'Sin especificar'
$ws.Range("I34").Value = 'Primera'
$ws.Range("J34").Value = 40
$ws.Range("K34").Value = 20000
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = 20000
$ws.Range("N34").Value = '$/caja 15 kilos granel'
$ws.Range("O34").Value = 'Región de Arica y Parinacota'
$ws.Range("P34").Value = 1333
$ws.Range("Q34").Value = 15
$ws.Range("R34").Value = 'Hortaliza'

# Row 35
$ws.Range("A35").Value = 10
$ws.Range("B35").Value = 'Vega Modelo de Temuco'
$ws.Range("C35").Value = 'La Araucanía'
$ws.Range("D35").Value = 44578
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = 100114002
$ws.Range("G35").Value = 'Camote'
$ws.Range("H35").Value = 'Sin especificar'
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 50
$ws.Range("K35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("M35").Value = 20000
$ws.Range("N35").Value = '$/malla 20 kilos'
$ws.Range("O35").Value = 'Perú'
$ws.Range("P35").Value = 1000
$ws.Range("Q35").Value = 20
$ws.Range("R35").Value = 'Hortaliza'

# Row 36
$ws.Range("A36").Value = 10
$ws.Range("B36").Value = 'Vega Modelo de Temuco'
$ws.Range("C36").Value = 'La Araucanía'
$ws.Range("D36").Value = 44315
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = 100114002
$ws.Range("G36").Value = 'Camote'
$ws.Range("H36").Value = 'Sin especificar'
$ws.Range("I36").Value = 'Primera'
$ws.Range("J36").Value = 30
$ws.Range("K36").Value = 20000
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = 20000
$ws.Range("N36").Value = '$/caja 15 kilos granel'
$ws.Range("O36").Value = 'Región de Arica y Parinacota'
$ws.Range("P36").Value = 1333
$ws.Range("Q36").Value = 15
$ws.Range("R36").Value = 'Hortaliza'

# Row 37
$ws.Range("A37").Value = 10
$ws.Range("B37").Value = 'Vega Modelo de Temuco'
$ws.Range("C37").Value = 'La Araucanía'
$ws.Range("D37").Value = 44315
$ws.Range("E37").Value = 9
$ws.Range("F37").Value = 100114002
$ws.Range("G37").Value = 'Camote'
$ws.Range("H37").Value = 'Sin especificar'
$ws.Range("I37").Value = 'Primera'
$ws.Range("J37").Value = 30
$ws.Range("K37").Value = 20000
$ws.Range("L37").Value = 20000
$ws.Range("M37").Value = 20000
$ws.Range("N37").Value = '$/malla 20 kilos'
$ws.Range("O37").Value = 'Región de Arica y Parinacota'
$ws.Range("P37").Value = 1000
$ws.Range("Q37").Value = 20
$ws.Range("R37").Value = 'Hortaliza'

# Row 38
$ws.Range("A38").Value = 10
$ws.Range("B38").Value = 'Vega Modelo de Temuco'
$ws.Range("C38").Value = 'La Araucanía'
$ws.Range("D38").Value = 44511
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = 100114002
$ws.Range("G38").Value = 'Camote'
$ws.Range("H38").Value = 'Sin especificar'
$ws.Range("I38").Value = 'Primera'
$ws.Range("J38").Value = 50
$ws.Range("K38").Value = 20000
$ws.Range("L38").Value = 20000
$ws.Range("M38").Value = 20000
$ws.Range("N38").Value = '$/malla 20 kilos'
$ws.Range("O38").Value = 'Perú'
$ws.Range("P38").Value = 1000
$ws.Range("Q38").Value = 20
$ws.Range("R38").Value = 'Hortaliza'

# Row 39
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = 'Vega Modelo de Temuco'
$ws.Range("C39").Value = 'La Araucanía'
$ws.Range("D39").Value = 44567
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = 100114002
$ws.Range("G39").Value = 'Camote'
$ws.Range("H39").Value = 'Sin especificar'
$ws.Range("I39").Value = 'Primera'
$ws.Range("J39").Value = 50
$ws.Range("K39").Value = 18000
$ws.Range("L39").Value = 18000
$ws.Range("M39").Value = 18000
$ws.Range("N39").Value = '$/malla 20 kilos'
$ws.Range("O39").Value = 'Región de Arica y Parinacota'
$ws.Range("P39").Value = 900
$ws.Range("Q39").Value = 20
$ws.Range("R39").Value = 'Hortaliza'

# Row 40
$ws.Range("A40").Value = 10
$ws.Range("B40").Value = 'Vega Modelo de Temuco'
$ws.Range("C40").Value = 'La Araucanía'
$ws.Range("D40").Value = 44455
$ws.Range("E40").Value = 9
$ws.Range("F40").Value = 100114002
$ws.Range("G40").Value = 'Camote'
$ws.Range("H40").Value = 'Sin especificar'
$ws.Range("I40").Value = 'Primera'
$ws.Range("J40").Value = 30
$ws.Range("K40").Value = 20000
$ws.Range("L40").Value = 20000
$ws.Range("M40").Value = 20000
$ws.Range("N40").Value = '$/malla 20 kilos'
$ws.Range("O40").Value = 'Perú'
$ws.Range("P40").Value = 1000
$ws.Range("Q40").Value = 20
$ws.Range("R40").Value = 'Hortaliza'

# Row 41
$ws.Range("A41").Value = 10
$ws.Range("B41").Value = 'Vega Modelo de Temuco'
$ws.Range("C41").Value = 'La Araucanía'
$ws.Range("D41").Value = 44329
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = 100114002
$ws.Range("G41").Value = 'Camote'
$ws.Range("H41").Value = 'Sin especificar'
$ws.Range("I41").Value = 'Primera'
$ws.Range("J41").Value = 40
$ws.Range("K41").Value = 20000
$ws.Range("L41").Value = 20000
$ws.Range("M41").Value = 20000
$ws.Range("N41").Value = '$/caja 15 kilos granel'
$ws.Range("O41").Value = 'Perú'
$ws.Range("P41").Value = 1333
$ws.Range("Q41").Value = 15
$ws.Range("R41").Value = 'Hortaliza'

# Row 42
$ws.Range("A42").Value = 10
$ws.Range("B42").Value = 'Vega Modelo de Temuco'
$ws.Range("C42").Value = 'La Araucanía'
$ws.Range("D42").Value = 44294
$ws.Range("E42").Value = 9
$ws.Range("F42").Value = 100114002
$ws.Range("G42").Value = 'Camote'
$ws.Range("H42").Value = 'Sin especificar'
$ws.Range("I42").Value = 'Primera'
$ws.Range("J42").Value = 5
$ws.Range("K42").Value = 20000
$ws.Range("L42").Value = 20000
$ws.Range("M42").Value = 20000
$ws.Range("N42").Value = '$/caja 15 kilos granel'
$ws.Range("O42").Value = 'Perú'
$ws.Range("P42").Value = 1333
$ws.Range("Q42").Value = 15
$ws.Range("R42").Value = 'Hortaliza'

# Row 43
$ws.Range("A43").Value = 10
$ws.Range("B43").Value = 'Vega Modelo de Temuco'
$ws.Range("C43").Value = 'La Araucanía'
$ws.Range("D43").Value = 44385
$ws.Range("E43").Value = 9
$ws.Range("F43").Value = 100114002
$ws.Range("G43").Value = 'Camote'
$ws.Range("H43").Value = 'Sin especificar'
$ws.Range("I43").Value = 'Primera'
$ws.Range("J43").Value = 18
$ws.Range("K43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("M43").Value = 20000
$ws.Range("N43").Value = '$/malla 20 kilos'
$ws.Range("O43").Value = 'Región de Arica y Parinacota'
$ws.Range("P43").Value = 1000
$ws.Range("Q43").Value = 20
$ws.Range("R43").Value = 'Hortaliza'

# Row 44
$ws.Range("A44").Value = 10
$ws.Range("B44").Value = 'Vega Modelo de Temuco'
$ws.Range("C44").Value = 'La Araucanía'
$ws.Range("D44").Value = 44498
$ws.Range("E44").Value = 9
$ws.Range("F44").Value = 100114002
$ws.Range("G44").Value = 'Camote'
$ws.Range("H44").Value = 'Sin especificar'
$ws.Range("I44").Value = 'Primera'
$ws.Range("J44").Value = 20
$ws.Range("K44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("M44").Value = 20000
$ws.Range("N44").Value = '$/malla 20 kilos'
$ws.Range("O44").Value = 'Región de Arica y Parinacota'
$ws.Range("P44").Value = 1000
$ws.Range("Q44").Value = 20
$ws.Range("R44").Value = 'Hortaliza'

# Row 45
$ws.Range("A45").Value = 10
$ws.Range("B45").Value = 'Vega Modelo de Temuco'
$ws.Range("C45").Value = 'La Araucanía'
$ws.Range("D45").Value = 44179
$ws.Range("E45").Value = 9
$ws.Range("F45").Value = 100114002
$ws.Range("G45").Value = 'Camote'
$ws.Range("H45").Value = 'Sin especificar'
$ws.Range("I45").Value = 'Primera'
$ws.Range("J45").Value = 20
$ws.Range("K45").Value = 20000
$ws.Range("L45").Value = 20000
$ws.Range("M45").Value = 20000
$ws.Range("N45").Value = '$/caja 15 kilos granel'
$ws.Range("O45").Value = 'Región de Arica y Parinacota'
$ws.Range("P45").Value = 1333
$ws.Range("Q45").Value = 15
$ws.Range("R45").Value = 'Hortaliza'

# Row 46
$ws.Range("A46").Value = 10
$ws.Range("B46").Value = 'Vega Modelo de Temuco'
$ws.Range("C46").Value = 'La Araucanía'
$ws.Range("D46").Value = 44595
$ws.Range("E46").Value = 9
$ws.Range("F46").Value = 100114002
$ws.Range("G46").Value = 'Camote'
$ws.Range("H46").Value = 'Sin especificar'
$ws.Range("I46").Value = 'Primera'
$ws.Range("J46").Value = 50
$ws.Range("K46").Value = 18000
$ws.Range("L46").Value = 18000
$ws.Range("M46").Value = 18000
$ws.Range("N46").Value = '$/malla 20 kilos'
$ws.Range("O46").Value = 'Perú'
$ws.Range("P46").Value = 900
$ws.Range("Q46").Value = 20
$ws.Range("R46").Value = 'Hortaliza'

# Ensure the date column keeps the existing date number format (style index 2),
# matching the rest of column D -- needed especially for the newly added row 46.
$ws.Range("D30:D46").NumberFormat = "YYYY-MM-DD HH:MM:SS"
